$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 515; this shifts the existing rows
# 515-548 down to 516-549 (values/formatting carried along automatically).
$ws.Rows.Item(515).Insert()

# Populate the newly inserted row 515 with the new data record.
$ws.Range("A515").Value = 10
$ws.Range("B515").Value = "Vega Modelo de Temuco"
$ws.Range("C515").Value = "La Araucanía"
$ws.Range("D515").Value = 44516
$ws.Range("E515").Value = 9
$ws.Range("F515").Value = "Fruta"
$ws.Range("G515").Value = 100109
$ws.Range("H515").Value = "Uva"
$ws.Range("I515").Value = 100109001
$ws.Range("J515").Value = "Uva"
$ws.Range("K515").Value = "Crimpson Seedless"
$ws.Range("L515").Value = "Primera"
$ws.Range("M515").Value = 55
$ws.Range("N515").Value = 20000
$ws.Range("O515").Value = 20000
$ws.Range("P515").Value = 20000
$ws.Range("Q515").Value = "$/bandeja 10 kilos"
$ws.Range("R515").Value = "Provincia de Limarí"
$ws.Range("S515").Value = 2000
$ws.Range("T515").Value = 10
